# Update the "想去人数" (number of interested attendees) figures on the
# "展览" and "全部类型" sheets to reflect the latest scraped data.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 1817
    $ws.Range("F3").Value = 8226
    $ws.Range("F5").Value = 315
}
